# "add delete ad by admin" - mark the Admin section (Home Screen / Approve /
# Reject / Edit / Delete Ad rows) as implemented by flagging column E
# ("Yes/No") with "YES" for rows 34-38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E34:E38").Value = "YES"

# Reflect the author's final cursor position / scroll (row 33 area -> row 45
# area, last touched cell E38) as recorded in the saved view state.
$win = $excel.Application.ActiveWindow
$win.ScrollRow = 30
$win.ScrollColumn = 1
$ws.Range("E38").Select()
